# Update odds values for rows 2-5 (matches/games) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Banfield vs Tigre)
$ws.Range("I2").Value  = 2.9
$ws.Range("AC2").Value = 7
$ws.Range("AE2").Value = 17
$ws.Range("AG2").Value = 7.5
$ws.Range("AJ2").Value = 29
$ws.Range("AN2").Value = 4.5
$ws.Range("AS2").Value = 251
$ws.Range("AX2").Value = 17

# Row 3 (Defensa y Justicia vs Dep. Riestra)
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.67

# Row 4 (Platense vs Godoy Cruz)
$ws.Range("G4").Value  = 2
$ws.Range("H4").Value  = 2.88
$ws.Range("J4").Value  = 2.88
$ws.Range("M4").Value  = 1.17
$ws.Range("N4").Value  = 4.75
$ws.Range("O4").Value  = 1.73
$ws.Range("P4").Value  = 2
$ws.Range("X4").Value  = 7.5
$ws.Range("AC4").Value = 4.75
$ws.Range("AD4").Value = 6
$ws.Range("AG4").Value = 8.5
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 13
$ws.Range("AW4").Value = 6

# Row 5 (Atl. Tucuman vs Huracan)
$ws.Range("L5").Value  = 4.33
$ws.Range("M5").Value  = 1.11
$ws.Range("N5").Value  = 6.5
$ws.Range("O5").Value  = 1.57
$ws.Range("P5").Value  = 2.25
$ws.Range("Q5").Value  = 2.7
$ws.Range("R5").Value  = 1.44
$ws.Range("S5").Value  = 1.62
$ws.Range("T5").Value  = 2.2
$ws.Range("AC5").Value = 6
$ws.Range("AG5").Value = 7.5
$ws.Range("AK5").Value = 34
$ws.Range("AP5").Value = 29
$ws.Range("AR5").Value = 81
$ws.Range("AT5").Value = 2.2
$ws.Range("AU5").Value = 9.5
$ws.Range("AX5").Value = 21
$ws.Range("BA5").Value = 126
$ws.Range("BB5").Value = 401
